# "Updates based on feedback"
#
# Row 34 (Body / search / word-paragraph-get-word-count / ...) had its
# MethodNameInTheSnippet value "searchAndCount" replaced with "run", and
# row 35's "parseText" was likewise replaced with "run". Both rows also
# lose the one-off "applyNumberFormat" cell style they had picked up
# (s="3"), reverting to the sheet's default/unstyled cells. The active
# selection moves from D35 to D34 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# Row 34: MethodNameInTheSnippet "searchAndCount" -> "run"
$ws.Range("D34").Value = "run"

# Row 35: MethodNameInTheSnippet "parseText" -> "run"
$ws.Range("D35").Value = "run"

# Both rows drop the stray numeric-format style they carried (s="3"),
# going back to the workbook's default cell style.
$ws.Range("A34:D35").Style = "Normal"

# Reflect the new active cell/selection (previously D35).
$ws.Range("D34").Select()
